$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Rows.Item(11).Insert()
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$ws1.Cells.Item(11,1).Value = "Jurisdiction"
$ws1.Cells.Item(11,2).Value = ""
